# Generate Report for Handoff
# Updates the status / priority / timestamp for the
# d798e671-f3e2-4daa-940d-943da01f2af3 record now that it is ready for
# handoff, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the d798e671-f3e2-4daa-940d-943da01f2af3.md record
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-12 06:16:46"

# --- zh-cn sheet: row 3 is the d798e671-f3e2-4daa-940d-943da01f2af3.md record
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-12 06:16:40"

# --- de-de sheet: row 3 is the d798e671-f3e2-4daa-940d-943da01f2af3.md record
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-12 06:16:46"

# --- Column width adjustments so the longer "Ready for handoff" status fits
$overview.Range("E1:F1").ColumnWidth = 16.3
$zhcn.Range("C1").ColumnWidth = 16.3
$dede.Range("C1").ColumnWidth = 16.3
